$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2866725.75
$ws.Range("C9").Value = 451971.99
$ws.Range("D9").Value = 3318697.74
$ws.Range("E9").Value = 13.61895615115584
$ws.Range("F9").Value = 86.38104384884416
$ws.Range("G9").Value = -56.31916080273851
$ws.Range("H9").Value = -48.23089545689565
$ws.Range("I9").Value = 28500
$ws.Range("J9").Value = 1220
$ws.Range("K9").Value = 29720
$ws.Range("L9").Value = 20504
$ws.Range("M9").Value = 161.85611295357
$ws.Range("N9").Value = 10.50246495491423
